$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.173706412315369
$ws.Range("B1").Value = 2.068058013916016
$ws.Range("C1").Value = 5.607265472412109
$ws.Range("D1").Value = 0.7988770604133606
$ws.Range("E1").Value = 0.9644936919212341
